$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("AQ3").Value = 'Die Daten der kompletten Zeitreihe wurden korrigiert.'
$ws.Range("AR3").Value = 'Data of the entire time series has been revised.'
$ws.Range("AT3").Value = '1990 zum Teil unsichere Datenbasis.<br>2020 vorläufige Daten.'
$ws.Range("AU3").Value = '1990 partly uncertain data basis.<br>2020 provisional data.'
$ws.Range("AW3").Value = '1992 zum Teil unsichere Datenbasis.<br>2018 vorläufige Daten.'
$ws.Range("AX3").Value = '1992 partly uncertain data basis.<br>2018 provisional data.'
$ws.Range("AQ6").Value = 'Todesfälle pro 100 000 Einwohner/ -innen unter 70 Jahren (ohne unter 1-Jährige).<br>Altersstandardisierung: alte Europastandardbevölkerung.'
$ws.Range("AR6").Value = 'Deaths per 100,000 inhabitants below the age of 70 (excluding those less than one year old).<br>Age-standardisation: old European standard population.'
$ws.Range("AQ7").Value = 'Raucherquote von Erwachsenen: 2021 vorläufige Daten.<br>Jugendliche: 12 bis 17 Jahre.<br> Erwachsene: ab 15 Jahren.'
$ws.Range("AR7").Value = 'Smoking rate among adults: 2021 preliminary data.<br>Adolescents: 12- to 17-year-olds.<br>Adults: 15 years and older.'
$ws.Range("AQ8").Value = 'Die Daten basieren auf einer Sonderauswertungund sind nicht öffentlich zugänglich.<br>Altersstandardisierung: Bevölkerungsstand: 31. Dezember 2015.'
$ws.Range("AR8").Value = 'Data is based on a special evaluation and are not publicly available.<br>Age-standardization: population status: December 31, 2015.'
$ws.Range("AQ9").Value = '2021 vorläufige Daten.<br>Altersstandardisierung: neue Europabevölkerung.'
$ws.Range("AR9").Value = '2021 preliminary data.<br>Age-standardization: new European standard population.'
$ws.Range("AQ10").Value = 'Gemittelter Index der Messzahlen: Schwefeldioxid (SO2), Stickstoffoxide (NOx), Ammoniak (NH3), flüchtige organische Verbindungen (NMVOC) und Feinstaub (PM2.5).'
$ws.Range("AR10").Value = 'Average index of measured values: Sulphur dioxide (SO2), nitrogen oxides (NOx), non-methane volatile organic compounds (NMVOCs) and particulate matter (PM2.5).'
$ws.Range("AQ17").Value = 'Frauen in Aufsichtsräten: Stand Januar des jeweiligen Jahres.<br>Frauen im öffentlichen Dienst des Bundes: 
- Die Daten basieren auf einer Sonderauswertung und sind nicht öffentlich zugänglich.
- Stand 30.06. des jeweiligen Jahres.
- Gleichberechtigte Teilhabe: Annähernd numerische Gleichheit.'
$ws.Range("AR17").Value = 'Women on supervisory boards: figures as at January each year.<br>Women in management positions in the federal civil service: 
- Data is based on a special evaluation and is not publicly available.
- Figures as at 30 June each year.
- Equal participation: approximate numerical equality.'
$ws.Range("AQ21").Value = 'Quelldaten auf Basis von Daten der Bund/Länder-Arbeitsgemeinschaft Wasser.<br>Basis EUA-Messnetz: Schwellenwert 50 Milligramm Nitrat pro Liter im Jahresmittel.<br> Berlin, Bremen und Hamburg (Stadtstaaten): keine Daten vorhanden. Grund: zu wenig Messstellen.'
$ws.Range("AR21").Value = 'Sourcedata based on data from the German Working Group
on Water Issues of the Länder and the Federal Government.<br>Basis EEA monitoring network: the threshold is an annual average of 50 mg nitrate per litre of groundwater.<br>Berlin, Bremen and Hamburg (city states): no data available. Reason: too few monitoring points.'
$ws.Range("AQ22").Value = 'Die Daten basieren auf einer Sonderauswertung und sind nicht öffentlich zugänglich.<br>Aufgrund methodischer Änderungen sind die Ergebnisse ab 2019 nur eingeschränkt mit den Vorjahren vergleichbar (Zeitreihenbruch).<br>Ab 2019 wird die Zeitreihe für die Bereiche Trinkwasserversorgung und Sanitärversorgung getrennt ausgewiesen.'
$ws.Range("AR22").Value = 'Data is based on a special evaluation and is not publicly available.<br>Due to methodological changes, the results from 2019 are only comparable with previous years to a limited extent.<br>Beginning in 2019, the time series for drinking water supply and sanitation will be reported separately.'
$ws.Range("AT23").Value = '2019 bis 2021 vorläufige Daten.'
$ws.Range("AU23").Value = '2019 to 2021 preliminary data.'
$ws.Range("AW23").Value = '2021 vorläufige Daten.'
$ws.Range("AX23").Value = '2021 preliminary data.'
$ws.Range("AQ26").Value = 'Gesamtrohstoffproduktivität:
- Das Ziel entspricht einer Beibehaltung des Trends der Jahre 2000 - 2010, welcher durchschnittlich rund 1,6 % Steigerung pro Jahr aufwies.
- Ab 2010 aufgrund methodischer Änderungen korrigierte Daten.<br>
Rohstoffeinsatz für Konsum, Investitionen und Exporte: Ab 2010 aufgrund methodischer Änderungen korrigierte Daten.<br>Wert von Konsum, Investitionen und Exporten (preisbereinigt): 2001 bis 2007 interpolierte Daten.'
$ws.Range("AR26").Value = 'Raw material input productivity:
- The target represents a continuation of the trend in the period from 2000 to 2010, when the average 
annual increase amounted to about 1.6%.
- From 2010 revised data due to methodological changes.<br>Raw material input for consumption, investment and exports: From 2010 revised data due to methodological changes.<br>Value of consumption, investment and exports (price-adjusted): 2001 to 2007 interpolated data.'
$ws.Range("AQ27").Value = 'Vorjahrespreise verkettet: 2015 = 100.<br>2019 bis 2021 vorläufige Daten.'
$ws.Range("AR27").Value = 'Previous years’ prices chain-linked: 2015 = 100.<br>Provisional data for 2019 to 2021.'
$ws.Range("AQ28").Value = '2019 bis 2021 vorläufige Daten.'
$ws.Range("AR28").Value = '2019 to 2021 provisional data.'
$ws.Range("AQ30").Value = 'Vorjahrespreise verkettet: 2015 = 100.<br>2019 bis 2021 vorläufige Daten.'
$ws.Range("AR30").Value = 'Previous years’ prices chain-linked: 2015 = 100.<br>2019 to 2021 provisional data.'
$ws.Range("AQ31").Value = ''
$ws.Range("AR31").Value = ''
$ws.Range("AT31").Value = 'Aufgrund einer umfassenden Neugestaltung des Mikrozensus ist ein Vergleich der Daten des Erhebungsjahres 2020 mit den Vorjahren nur eingeschränkt möglich (Zeitreihenbruch).'
$ws.Range("AU31").Value = 'Due to comprehensive redesign of the microcensus it is not possible to compare the data of the survey year 2020 with previous years (break in time series).'
$ws.Range("AW31").Value = 'Aufgrund einer umfassenden Neugestaltung des Mikrozensus ist ein Vergleich der Daten des Erhebungsjahres 2020 mit den Vorjahren nur eingeschränkt möglich (Zeitreihenbruch).<br>Bundesländer: Die Daten basieren auf einer Sonderauswertung und sind nicht öffentlich zugänglich.<br>Bremen: 
- Männer 2010 eingeschränkter Aussagewert.
- Frauen 2010 bis 2012 sowie 2015 eingeschränkter Aussagewert.<br>Saarland: 
- Frauen 2010 eingeschränkter Aussagewert.'
$ws.Range("AX31").Value = 'Due to comprehensive redesign of the microcensus it is not possible to compare the data of the survey year 2020 with previous years (break in time series).<br>Länder: The data is based on a special evaluation and is not publicly available.<br>Bremen: 
- Men 2010 limited significance.
- Women 2010 to 2012 and 2015 limited significance.<br>Saarland: 
- Women 2010 limited significance.'
$ws.Range("AQ34").Value = 'FTTB/H: Fibre-to-the-Building/Home (Glasfaser bis zum Gebäude/in die Wohnung).<br>CATV: Cable Television (Kabelfernsehen).'
$ws.Range("AQ36").Value = 'Die bislang separat durchgeführte Erhebung „Leben in Europa“ (EU-SILC) wurde 2020 in den Mikrozensus als Unterstichprobe integriert. Durch den Wechsel von einer freiwilligen zu einer in Teilen auskunftspflichtigen Befragung verbunden mit einer neuen Stichprobenzusammensetzung ist ein Vergleich der Daten des Erhebungsjahres 2020 mit den Vorjahren nicht möglich (Zeitreihenbruch).<br>Gini-Koeffizient des verfügbaren Äquivalenzeinkommens vor Sozialleistungen: Renten von den Sozialleistungen ausgeschlossen.<br>Gini-Koeffizient des verfügbaren Äquivalenzeinkommens (EU): 
- Für EU: 2019 von Eurostat geschätzte Daten. 
- Ab 2020: EU-27 (ohne Vereinigtes Königreich).'
$ws.Range("AR36").Value = 'The ''Leben in Europa'' survey (German name of the European Union Statistics on Income and Living Conditions - EU-SILC), which was conducted separately in the past, was integrated as a subsample into the microcensus in 2020. Comparing the data of reference year 2020 with those of previous years is not possible (break in the time series) as the voluntary survey was changed over to a partly compulsory survey and the composition of the sample was changed.<br>Gini coefficient of equivalised disposable income before social transfers: Pensions not included in social benefits.<br>Gini coefficient of equivalised disposable income:
- For EU: 2019 data estimated by Eurostat.
- From 2020: EU-27 (without UK).'
$ws.Range("AQ37").Value = 'Siedlungs- und Verkehrsfläche: Aufgrund methodischer Änderungen in der amtlichen Flächenerhebung (Einführung Amtliches Liegenschaftskataster-Informationssystem (ALKIS) ab Berichtsjahr 2016) sind die Ergebnisse ab 2016 nur eingeschränkt mit den Vorjahren vergleichbar (Zeitreihenbruch).'
$ws.Range("AR37").Value = 'Settlement and transport area: Due to methodological changes in the official survey of land (introduction of the official land register information system (ALKIS) from the 2016 reporting year), the results from 2016 are only comparable with previous years to a limited extent (break in the time series).'
$ws.Range("AQ38").Value = 'Siedlungs- und Verkehrsfläche: Aufgrund methodischer Änderungen in der amtlichen Flächenerhebung (Einführung Amtliches Liegenschaftskataster-Informationssystem (ALKIS) ab Berichtsjahr 2016) sind die Ergebnisse ab 2016 nur eingeschränkt mit den Vorjahren vergleichbar (Zeitreihenbruch).<br>Es handelt sich um den nach Einwohnerinnen und Einwohnern gewichteten gleitenden Vierjahresdurchschnitt der ländlichen und nicht ländlichen Räume.'
$ws.Range("AR38").Value = 'Settlement and transport area: Due to methodological changes in the official survey of land (introduction of the official land register information system (ALKIS) from the 2016 reporting year), the results from 2016 are only comparable with previous years to a limited extent (break in time series).<br>The series shows the four-year moving average of rural and non-rural areas weighted by population.'
$ws.Range("AQ44").Value = 'Die Daten basieren auf einer Sonderauswertung. Die Grunddaten hierfür sind öffentlich zugänglich.'
$ws.Range("AR44").Value = 'Data is based on a special evaluation. The basic data for this is publicly available.'
$ws.Range("AQ47").Value = 'EMAS: Eco-Management and Audit Scheme (Umweltmanagementsystem).'
$ws.Range("AQ48").Value = 'Die Daten basieren auf einer Sonderauswertung und sind nicht öffentlich zugänglich.'
$ws.Range("AW48").Value = 'Aufgrund methodischer Änderungen sind die Ergebnisse ab 2016 nur eingeschränkt mit dem Vorjahr vergleichbar.'
$ws.Range("AX48").Value = 'Due to methodological changes, the results from 2016 are only comparable with previous year to a limited extent.'
$ws.Range("AQ50").Value = 'Anpassung des Zieljahres und -wertes von vier Milliarden Euro bis 2020 gemäß Grundsatzbeschluss 2022.'
$ws.Range("AR50").Value = 'Adjustment of the target year and value of four billion euros by 2020 in accordance with the 2022 policy decision.'
$ws.Range("AQ52").Value = 'MSY: Maximum Sustainable Yield (höchstmöglicher Dauerertrag).'
$ws.Range("AQ57").Value = 'Die Daten basieren auf einer Sonderauswertung. Die Grunddaten hierfür sind öffentlich zugänglich.'
$ws.Range("AR57").Value = 'Data is based on a special evaluation. The basic data for this is publicly available.'
$ws.Range("AT58").Value = 'Partnerländer einschließlich Südsudan.<br>Länderliste für die bilaterale staatliche Entwicklungszusammenarbeit des BMZ wurde 2020 aktualisiert. Es werden 10 Länder mehr aufgeführt zuzüglich China und Sudan.'
$ws.Range("AQ60").Value = 'Das Ziel entspricht einer Steigerung der Anzahl an Studierenden und Forschenden um 10 % gegenüber 2015 in 2020.<br>Die Daten zu den Forschenden basieren auf einer Sonderauswertung und sind nicht öffentlich zugänglich.<br>LDC: Least Developed Countries (am wenigsten entwickelte Länder).'
$ws.Range("AR60").Value = 'The target for 2020 corresponds to a 10% increase in the number of students and researchers compared to 2015.<br>Data for researchers is based on a special evaluation and not publicly availabe.<br>LDC: Least Developed Countries.'
$ws.Range("AQ61").Value = 'Die Daten basieren auf einer Sonderauswertung. Die Grunddaten hierfür sind öffentlich zugänglich.<br>2021 vorläufige Daten.<br>LDC: Least Developed Countries (am wenigsten entwickelte Länder).'
$ws.Range("AR61").Value = 'Data is based on a special evaluation. The basic data for this is publicly available.<br>2021 provisional data. <br>LDC: least developed countries.'

$ws.Range("K34").Copy()
$ws.Range("BN34").PasteSpecial(-4122)
$ws.Range("BN34").Value = 10
